$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column C ("Post Treatment" bucket for session 3 / pre-treatment phase
# self-report measurements) for rows 2-15.
$values = @{
    2  = "A little worse"
    3  = "Somewhat worse"
    4  = "A little worse"
    5  = "A little worse"
    6  = "Somewhat worse"
    7  = "A lot worse"
    8  = "Not worse"
    9  = "Somewhat worse"
    10 = "Somewhat worse"
    11 = "Not worse"
    12 = "Not worse"
    13 = "A lot worse"
    14 = "A little worse"
    15 = "Not worse"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}

# Move the active selection to C16, matching the saved cursor position.
$ws.Range("C16").Select()
